$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 values ---
$ws.Range("B2").Value = "Tac1"
$ws.Range("C2").Value = "Tacr3"
$ws.Range("D2").Value = "MuSCs"

$ws.Range("I2").Value = 0.9476581720434079
$ws.Range("J2").Value = 0.947658172043408
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.0556235
$ws.Range("N2").Value = 0.111247
$ws.Range("Q2").Value = 0.1379606679453333
$ws.Range("R2").Value = 0.827764007672
$ws.Range("S2").Value = 0.9476581720434079
$ws.Range("T2").Value = 0.947658172043408

# --- Add new row 3 ---
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("B3").Value = "Tac1"
$ws.Range("C3").Value = "Tacr3"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1369916666666667
$ws.Range("H3").Value = 0.410975
$ws.Range("I3").Value = 0.0523418279565921
$ws.Range("J3").Value = 0.0523418279565921
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.0556235
$ws.Range("N3").Value = 0.111247
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.007619955970833332
$ws.Range("R3").Value = 0.045719735825
$ws.Range("S3").Value = 0.0523418279565921
$ws.Range("T3").Value = 0.0523418279565921
